$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old hyperlinks bound to the old layout (D2,H2,D3,H3) ---
$ws.Range("D2").Hyperlinks.Delete()
$ws.Range("H2").Hyperlinks.Delete()
$ws.Range("D3").Hyperlinks.Delete()
$ws.Range("H3").Hyperlinks.Delete()

# --- Wipe all existing cell content/styling so we can rebuild the new layout cleanly ---
$ws.Cells.Clear()

# --- Pre-format as Text ("@") before writing any values so numeric-looking ---
# --- strings (ids, phone numbers) are preserved as text, not auto-converted ---
# --- to numbers. E2 and C3:F3 stay General since they hold real numbers, and ---
# --- A3 is intentionally left untouched to match the template (it keeps the ---
# --- default/general style with no explicit style index). ---
$ws.Range("A1:I1").NumberFormat = "@"
$ws.Range("A2,B2,C2,D2,F2,G2,H2,I2").NumberFormat = "@"
$ws.Range("B3,G3,H3,I3").NumberFormat = "@"

# --- Header row ---
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "email"
$ws.Range("C1").Value = "work_phone"
$ws.Range("D1").Value = "mobile_phone"
$ws.Range("E1").Value = "direct_line"
$ws.Range("F1").Value = "fax"
$ws.Range("G1").Value = "website"
$ws.Range("H1").Value = "title"
$ws.Range("I1").Value = "staff_id"

# --- Row 2 ---
$ws.Range("A2").Value = "James Y.T. Hu"
$ws.Range("B2").Value = "jameshu@pyengineering.com"
$ws.Range("C2").Value = "28383030"
$ws.Range("D2").Value = "00000000"
$ws.Range("E2").Value = 12345600
$ws.Range("F2").Value = "12345678"
$ws.Range("G2").Value = "www.pyengineering.com"
$ws.Range("H2").Value = "Sample Guy"
$ws.Range("I2").Value = "00001"

# --- Row 3 ---
$ws.Range("A3").Value = "John N.B. Doe"
$ws.Range("B3").Value = "anotherone@pyengineering.com"
$ws.Range("C3").Value = 28282828
$ws.Range("D3").Value = 65656565
$ws.Range("E3").Value = 12345600
$ws.Range("F3").Value = 12348765
$ws.Range("G3").Value = "www.pyengineering.com"
$ws.Range("H3").Value = "Sample Guy"
$ws.Range("I3").Value = "00002"

# --- Hyperlinks on the new layout ---
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:jameshu@pyengineering.com")
$ws.Hyperlinks.Add($ws.Range("G2"), "http://www.pyengineering.com/")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:anotherone@pyengineering.com")
$ws.Hyperlinks.Add($ws.Range("G3"), "http://www.pyengineering.com/")

# --- Column A width so the full names are visible ---
$ws.Columns.Item(1).ColumnWidth = 20.6

# --- Match the saved selection in the template ---
$ws.Range("L7").Select()
